# Add a new "2022-Q3" sheet (fund holdings for the new quarter), inserted
# right after "总计" / before "2022-Q2", and record the new quarter's
# summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Several columns on the per-quarter sheets store numeric-looking
    # figures as text (t="inlineStr" in the original file). Range.Value
    # auto-detects numeric-looking strings as numbers, so force the Text
    # number format first, then drop back to the Normal style so we don't
    # leave a stray "@" format marker behind.
    $rng.NumberFormat = "@"
    $rng.Value = [string]$val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by cloning "2022-Q2" (keeps headers,
#    column-A styling, borders, etc. identical) and inserting it right
#    before "2022-Q2" in the tab order.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2: 512660 国泰中证军工ETF
Set-TextValue $q3.Range("D2") "105.57"
Set-TextValue $q3.Range("E2") "99.74"
Set-TextValue $q3.Range("F2") "4.04"
Set-TextValue $q3.Range("G2") "4.2650"
$q3.Range("H2").Value = 6

# Row 3: 161024 富国中证军工指数A
$q3.Range("C3").Value = "富国中证军工指数A"
Set-TextValue $q3.Range("D3") "51.00"
Set-TextValue $q3.Range("E3") "93.67"
Set-TextValue $q3.Range("F3") "3.82"
Set-TextValue $q3.Range("G3") "1.9482"
$q3.Range("H3").Value = 6

# Row 4: 512680 广发中证军工ETF
Set-TextValue $q3.Range("D4") "30.69"
Set-TextValue $q3.Range("E4") "99.32"
Set-TextValue $q3.Range("F4") "4.04"
Set-TextValue $q3.Range("G4") "1.2399"
$q3.Range("H4").Value = 6

# Row 5: 163115 申万菱信中证军工指数A
Set-TextValue $q3.Range("D5") "11.85"
Set-TextValue $q3.Range("E5") "93.33"
Set-TextValue $q3.Range("F5") "3.81"
Set-TextValue $q3.Range("G5") "0.4515"
$q3.Range("H5").Value = 6

# Row 6: 501019 国泰国证航天军工指数（LOF）A
$q3.Range("C6").Value = "国泰国证航天军工指数（LOF）A"
Set-TextValue $q3.Range("D6") "9.07"
Set-TextValue $q3.Range("E6") "93.51"
Set-TextValue $q3.Range("F6") "3.95"
Set-TextValue $q3.Range("G6") "0.3583"
$q3.Range("H6").Value = 7

# Row 7: 502003 易方达军工指数（LOF）A
Set-TextValue $q3.Range("D7") "6.92"
Set-TextValue $q3.Range("E7") "94.39"
Set-TextValue $q3.Range("F7") "3.85"
Set-TextValue $q3.Range("G7") "0.2664"
$q3.Range("H7").Value = 6

# Row 8: 512560 易方达中证军工ETF
Set-TextValue $q3.Range("D8") "5.97"
Set-TextValue $q3.Range("E8") "98.91"
Set-TextValue $q3.Range("F8") "4.02"
Set-TextValue $q3.Range("G8") "0.2400"
$q3.Range("H8").Value = 6

# Row 9: 512810 华宝中证军工ETF (renamed from 华宝兴业中证军工ETF)
$q3.Range("C9").Value = "华宝中证军工ETF"
Set-TextValue $q3.Range("D9") "3.92"
Set-TextValue $q3.Range("E9") "98.39"
Set-TextValue $q3.Range("F9") "4.00"
Set-TextValue $q3.Range("G9") "0.1568"
$q3.Range("H9").Value = 6

# Row 10: 013035 富国中证军工指数C
Set-TextValue $q3.Range("D10") "1.16"
Set-TextValue $q3.Range("E10") "93.67"
Set-TextValue $q3.Range("F10") "3.82"
Set-TextValue $q3.Range("G10") "0.0443"
$q3.Range("H10").Value = 6

# Row 11: 012842 易方达军工指数（LOF）C
Set-TextValue $q3.Range("D11") "1.04"
Set-TextValue $q3.Range("E11") "94.39"
Set-TextValue $q3.Range("F11") "3.85"
Set-TextValue $q3.Range("G11") "0.0400"
$q3.Range("H11").Value = 6

# Row 12: now 015599 国泰国证航天军工指数（LOF）C (swapped with old row 13)
$q3.Range("B12").Value = "015599"
$q3.Range("C12").Value = "国泰国证航天军工指数（LOF）C"
Set-TextValue $q3.Range("D12") "0.03"
Set-TextValue $q3.Range("E12") "93.51"
Set-TextValue $q3.Range("F12") "3.95"
Set-TextValue $q3.Range("G12") "0.0012"
$q3.Range("H12").Value = 7

# Row 13: now 016209 申万菱信中证军工指数C (swapped with old row 12)
$q3.Range("B13").Value = "016209"
$q3.Range("C13").Value = "申万菱信中证军工指数C"
Set-TextValue $q3.Range("D13") "0.01"
Set-TextValue $q3.Range("E13") "93.33"
Set-TextValue $q3.Range("F13") "3.81"
Set-TextValue $q3.Range("G13") "0.0004"
$q3.Range("H13").Value = 6

# ---------------------------------------------------------------------
# 2) Insert the new totals row into "总计" (sheet 1), above the existing
#    "2022-Q2" row, and fill in the 2022-Q3 figures.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 9.01

# ---------------------------------------------------------------------
# 3) Restore the originally-active tab ("2020-Q4", the last sheet) so the
#    freshly-inserted "2022-Q3" sheet doesn't steal tab focus.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
